$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1.01
$ws.Range("J2").Value = 3.6
$ws.Range("L2").Value = 1.3
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 1.83
$ws.Range("O2").Value = 1.01
$ws.Range("R2").Value = 1.18
$ws.Range("S2").Value = 1.01
$ws.Range("T2").Value = 1.01
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 1.03
$ws.Range("W2").Value = 2.26
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000
$ws.Range("L3").Value = 1.34
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.66
$ws.Range("O3").Value = 1.01
$ws.Range("R3").Value = 1.18
$ws.Range("S3").Value = 1.89
$ws.Range("T3").Value = 1.01
$ws.Range("U3").Value = 1.01
$ws.Range("V3").Value = 1.33
$ws.Range("W3").Value = 1.5
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000
$ws.Range("G4").Value = 3.4
$ws.Range("I4").Value = 3.9
$ws.Range("S4").Value = 3.6
$ws.Range("W4").Value = 1.42
$ws.Range("F5").Value = 3.25
$ws.Range("H5").Value = 1.65
$ws.Range("W5").Value = 1.29
$ws.Range("I6").Value = 2.48
$ws.Range("P6").Value = 1.8
$ws.Range("Q6").Value = 1.86
$ws.Range("H7").Value = 3.8
$ws.Range("P7").Value = 1.75
$ws.Range("U8").Value = 1.91
$ws.Range("AA8").Value = 34
$ws.Range("AG8").Value = 18
$ws.Range("AH8").Value = 22
$ws.Range("AO8").Value = 28
$ws.Range("J9").Value = 3.7
$ws.Range("K9").Value = 4.1
$ws.Range("P9").Value = 2.2
$ws.Range("Q9").Value = 1.01
$ws.Range("F10").Value = 2.02
$ws.Range("I10").Value = 3.9
$ws.Range("F11").Value = 2.68
$ws.Range("H11").Value = 2.2
$ws.Range("K11").Value = 6
$ws.Range("P11").Value = 1.73
$ws.Range("Q11").Value = 1.84
$ws.Range("P12").Value = 2.06
$ws.Range("Q12").Value = 1.52
$ws.Range("F14").Value = 1.04
$ws.Range("P14").Value = 2.3
$ws.Range("Q14").Value = 1.43
$ws.Range("AJ18").Value = 29
$ws.Range("H19").Value = 2.6
$ws.Range("K19").Value = 3.9
$ws.Range("U20").Value = 2.2
$ws.Range("X20").Value = 18
$ws.Range("AD20").Value = 19.5
$ws.Range("P21").Value = 1.7
$ws.Range("H22").Value = 2.1
$ws.Range("K22").Value = 4.4
